$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. "1.00", "62.298.26") that must
# remain exact text, not be auto-coerced into numbers by Excel. Force the
# column to Text format before writing the values.
$ws.Range("D2:D51").NumberFormat = "@"

$data = @(
  ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "62.298.26", "  +2.32%  ")
  ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.413.32", "  -0.63%  ")
  ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.01", "  +0.87%  ")
  ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "572.15", "  +1.06%  ")
  ,@(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "144.59", "  +4.32%  ")
  ,@(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.998", "  -0.46%  ")
  ,@(8, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.538", "  +0.29%  ")
  ,@(9, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "2.438.25", "  +0.93%  ")
  ,@(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.110", "  +4.13%  ")
  ,@(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.161", "  +0.76%  ")
  ,@(12, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "5.26", "  +3.91%  ")
  ,@(13, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.350", "  +3.63%  ")
  ,@(14, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "26.81", "  +3.73%  ")
  ,@(15, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000179", "  +7.09%  ")
  ,@(16, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.853.94", "  +1.43%  ")
  ,@(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "61.999.91", "  +2.09%  ")
  ,@(18, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.436.39", "  +0.74%  ")
  ,@(19, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "7.95", "  -4.93%  ")
  ,@(20, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "10.86", "  +2.05%  ")
  ,@(21, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "326.83", "  +0.61%  ")
  ,@(22, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.13", "  +2.07%  ")
  ,@(23, "SuiNetwork", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui", "2.03", "  +12.07%  ")
  ,@(24, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  -0.22%  ")
  ,@(25, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "65.42", "  +1.37%  ")
  ,@(26, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "621.04", "  +10.62%  ")
  ,@(27, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "8.48", "  +5.00%  ")
  ,@(28, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0991", "  +7.66%  ")
  ,@(29, "Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "1.00", "  +0.24%  ")
  ,@(30, "WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "2.527.48", "  +0.04%  ")
  ,@(31, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "8.09", "  +2.06%  ")
  ,@(32, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "1.41", "  +7.83%  ")
  ,@(33, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.83", "  +1.65%  ")
  ,@(34, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.136", "  +3.23%  ")
  ,@(35, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.48", "  +3.66%  ")
  ,@(36, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.996", "  -0.82%  ")
  ,@(37, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "4.76", "  +4.76%  ")
  ,@(38, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "152.98", "  +0.83%  ")
  ,@(39, "PolygonEcosystemToken", "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol", "0.373", "  +1.02%  ")
  ,@(40, "RenderToken", "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render", "5.41", "  +5.39%  ")
  ,@(41, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "18.59", "  +1.58%  ")
  ,@(42, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "2.68", "  +15.21%  ")
  ,@(43, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.75", "  +5.36%  ")
  ,@(44, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "42.28", "  +0.61%  ")
  ,@(45, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "0.998", "  -0.08%  ")
  ,@(46, "BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.0₆0282", "  -2.83%  ")
  ,@(47, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "143.69", "  -0.67%  ")
  ,@(48, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.59", "  +2.27%  ")
  ,@(49, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "20.34", "  +6.37%  ")
  ,@(50, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.600", "  +2.07%  ")
  ,@(51, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0514", "  +3.00%  ")
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
}
